# Update row 9 (year 2025) of Sheet1 with refreshed "faturamento" figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value = 2955685.33
$ws.Range("C9").Value = 468266.99
$ws.Range("D9").Value = 3423952.32
$ws.Range("E9").Value = 13.67621234865794
$ws.Range("F9").Value = 86.32378765134206
$ws.Range("G9").Value = -54.7443303033543
$ws.Range("H9").Value = -46.62440840555122
$ws.Range("I9").Value = 29532
$ws.Range("J9").Value = 1261
$ws.Range("K9").Value = 30793
$ws.Range("L9").Value = 21246
$ws.Range("M9").Value = 161.1575035300763
$ws.Range("N9").Value = 10.02551007240733

$wb.Save()
